# Insert a new weekly price record as row 29 in the Cilantro price sheet,
# shifting the existing rows 29:65 down to 30:66 (dimension grows to A1:R66).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(29).Insert()

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44589
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 270
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 2500
$ws.Range("M29").Value = 2250
$ws.Range("N29").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 1125
$ws.Range("Q29").Value = 2
$ws.Range("R29").Value = "Hortaliza"
